$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "283.53"
Set-TextValue $ws.Range("E2") "1.90%"
Set-TextValue $ws.Range("D3") "28.29"
Set-TextValue $ws.Range("E3") "3.90%"
Set-TextValue $ws.Range("D4") "5.025"
Set-TextValue $ws.Range("E4") "3.45%"
Set-TextValue $ws.Range("D5") "0.06524"
Set-TextValue $ws.Range("E5") "1.85%"
Set-TextValue $ws.Range("D6") "7.239"
Set-TextValue $ws.Range("E6") "3.57%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D7") "1.431"
Set-TextValue $ws.Range("E7") "19.35%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D8") "0.9222"
Set-TextValue $ws.Range("E8") "4.94%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D9") "0.1545"
Set-TextValue $ws.Range("E9") "1.32%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.06560"
Set-TextValue $ws.Range("E10") "27.33%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.07469"
Set-TextValue $ws.Range("E11") "-0.45%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D12") "0.02771"
Set-TextValue $ws.Range("E12") "-5.67%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D13") "0.08972"
Set-TextValue $ws.Range("E13") "0.11%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D14") "0.001581"
Set-TextValue $ws.Range("E14") "0.32%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D15") "0.0006414"
Set-TextValue $ws.Range("E15") "0.34%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.006062"
Set-TextValue $ws.Range("E16") "-0.38%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.447"
Set-TextValue $ws.Range("E17") "-0.82%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D18") "3.378"
Set-TextValue $ws.Range("E18") "2.25%"
Set-TextValue $ws.Range("E19") "-0.33%"
Set-TextValue $ws.Range("E20") "3.43%"
Set-TextValue $ws.Range("D21") "0.1305"
Set-TextValue $ws.Range("E21") "-1.47%"
Set-TextValue $ws.Range("D22") "3.989"
Set-TextValue $ws.Range("E22") "2.23%"
Set-TextValue $ws.Range("D23") "0.1530"
Set-TextValue $ws.Range("E23") "1.70%"
Set-TextValue $ws.Range("D24") "0.04426"
Set-TextValue $ws.Range("E24") "0.66%"
Set-TextValue $ws.Range("D25") "0.001187"
Set-TextValue $ws.Range("E25") "1.05%"
Set-TextValue $ws.Range("D26") "0.004422"
Set-TextValue $ws.Range("E26") "13.51%"
Set-TextValue $ws.Range("D27") "0.0001252"
Set-TextValue $ws.Range("E27") "6.04%"
Set-TextValue $ws.Range("D28") "0.0001620"
Set-TextValue $ws.Range("E28") "-1.50%"
Set-TextValue $ws.Range("D40") "0.04126"
Set-TextValue $ws.Range("E40") "1.11%"
Set-TextValue $ws.Range("D41") "0.006741"
Set-TextValue $ws.Range("E41") "-0.85%"
Set-TextValue $ws.Range("D42") "0.1227"
Set-TextValue $ws.Range("E42") "4.53%"
Set-TextValue $ws.Range("D43") "0.002173"
Set-TextValue $ws.Range("E43") "15.01%"
Set-TextValue $ws.Range("E44") "2.38%"
Set-TextValue $ws.Range("D45") "0.00005627"
Set-TextValue $ws.Range("E45") "5.12%"
Set-TextValue $ws.Range("E46") "25.93%"
Set-TextValue $ws.Range("D47") "0.01301"
Set-TextValue $ws.Range("E47") "-29.72%"
